# household_new.xlsx edit
#
# Commit message: "Remove cell_type usage. Add expression sanity checking to
# XLSXconverter -- obviously bad formulas will be detected early."
#
# On the "survey" sheet the two "*.cell_type" helper columns (I:
# selectionArgs.cell_type, K: auxillaryHash.cell_type) are no longer used,
# so they are deleted outright (not just cleared) which shifts every column
# to their right one/two places to the left. Deleting real columns (rather
# than just clearing cell contents) is what reproduces the shrunk
# dimension/shared-string table seen in the diff.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item(1)   # survey
$choices = $wb.Worksheets.Item(2)   # choices
$settings = $wb.Worksheets.Item(3)  # settings (untouched)
$model   = $wb.Worksheets.Item(4)   # model

# Delete column K ("auxillaryHash.cell_type") first, then column I
# ("selectionArgs.cell_type"), so both removals land on the columns the
# diff shows disappearing (old I and old K).
$survey.Columns("K:K").Delete() | Out-Null
$survey.Columns("I:I").Delete() | Out-Null

# Reproduce the recorded cursor/selection state for each sheet.
$survey.Range("H11").Select() | Out-Null
$choices.Range("A1").Select() | Out-Null
$model.Range("C20").Select() | Out-Null

# "model" becomes the active tab/sheet (activeTab="3" in workbook.xml,
# tabSelected="1" on sheet4).
$model.Activate() | Out-Null
